# South Korea K3 League - base update (17-02-2024 22:47)
#
# The underlying data rows keep their sequential "id" (column A) fixed in
# place, but the rest of each row's fields (columns B..AC) were reshuffled
# among small clusters of adjacent rows (mostly simple pair swaps, plus two
# 3-row rotations). This script reproduces that reshuffle by reading the
# B:AC range of the involved rows into memory and writing it back out in the
# new arrangement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return $ws.Range("B$row`:AC$row").Value2
}

function Set-RowData($row, $data) {
    $ws.Range("B$row`:AC$row").Value2 = $data
}

function Swap-Rows($rowA, $rowB) {
    $dataA = Get-RowData $rowA
    $dataB = Get-RowData $rowB
    Set-RowData $rowA $dataB
    Set-RowData $rowB $dataA
}

# Simple pairwise swaps (id / column A stays put on each row).
$pairs = @(
    @(11, 12),
    @(14, 15),
    @(39, 40),
    @(50, 51),
    @(94, 95),
    @(110, 111),
    @(124, 125),
    @(132, 133),
    @(136, 137),
    @(176, 177),
    @(201, 202)
)

foreach ($p in $pairs) {
    Swap-Rows $p[0] $p[1]
}

# Three-way rotations: row 204 receives what was in row 206,
# row 205 receives what was in row 204, row 206 receives what was in row 205.
$d204 = Get-RowData 204
$d205 = Get-RowData 205
$d206 = Get-RowData 206
Set-RowData 204 $d206
Set-RowData 205 $d204
Set-RowData 206 $d205

# row 208 receives what was in row 210, row 209 receives what was in row 208,
# row 210 receives what was in row 209.
$d208 = Get-RowData 208
$d209 = Get-RowData 209
$d210 = Get-RowData 210
Set-RowData 208 $d210
Set-RowData 209 $d208
Set-RowData 210 $d209
